$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gear cards removed; their functionality is replaced by newly added skill cards,
# and the "Finger" monster card is renamed to "Sealed Finger".

# A18: card name "手指" (Finger) -> "封印手指" (Sealed Finger)
$ws.Range("A18").Value = "封印手指"

# F18: effect text now refers to the quoted card name "《封印手指》牌" instead of "手指牌"
$ws.Range("F18").Value = "持续：不能被重抽。<br>`n离场时：回到玩家手牌。<br>`n从手牌发动：如果手牌中有合计5张《封印手指》牌，移除当前挑战牌。"

# F15 (Collector): old Gear-era effect replaced with a new skill effect about item cards
$ws.Range("F15").Value = "持续：玩家不能使用道具牌。<br>`n进入墓地时：移除当前挑战牌。"

# F16 (Sandworm): "交锋中" (during the clash) changed to "同一行中" (same row)
$ws.Range("F16").Value = "交锋前：同一行中其他牌点数减1。<br>`n压制时：本牌点数减1，然后后退一行。"

$ws.Range("F17").Select()
